$d = $word.ActiveDocument

# --- Add the three new character styles -------------------------------
$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "2022 Fechas de la campaña..." run -------
$datesText = "2022 Fechas de la campaña para Constelación de Hércules: 13-22 de junio, 12-21 de julio, 10-19 de agosto"
$rng = $d.Content
$rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($rng.Find.Found) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $rng.End = $d.Content.End
    $rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNParagraph to the "Usted está participando..." run -------
$paragraphText = "Usted está participando en una campaña mundial para observar y registrar las estrellas visibles más débiles como un medio para medir la contaminación lumínica en un lugar determinado. Localizando y observando la  Constelación de Hércules en el cielo nocturno y comparándolo con las cartas estelares, la gente de todo el mundo aprenderán cómo las luces de su comunidad contribuyen a la contaminación lumínica. Sus contribuciones a la base de datos en línea documentarán el cielo nocturno visible."
$rng2 = $d.Content
$rng2.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($rng2.Find.Found) {
    $rng2.Style = "GaNParagraph"
    $rng2.Collapse(0)
    $rng2.End = $d.Content.End
    $rng2.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNLinks to the campaign map link run -----------------------
$linkText = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3 = $d.Content
$rng3.Find.Execute($linkText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($rng3.Find.Found) {
    $rng3.Style = "GaNLinks"
    $rng3.Collapse(0)
    $rng3.End = $d.Content.End
    $rng3.Find.Execute($linkText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

Write-Output "done"
